$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right under the page title.
# ------------------------------------------------------------------
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -ne -1) {
    $d.Paragraphs.Item($metaIndex).Range.Delete()
}

# ------------------------------------------------------------------
# 2. Find the trailing paragraph that holds the (italic) AI image
#    prompt - it is the very last paragraph in the document.
# ------------------------------------------------------------------
$imageIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Create a feature image")) {
        $imageIndex = $i
        break
    }
}

if ($imageIndex -ne -1) {
    $imagePara = $d.Paragraphs.Item($imageIndex)

    # ---------------------------------------------------------------
    # 2a. Swap the image-prompt text for the meta description text,
    #     keeping the paragraph's existing (italic) formatting intact.
    #     We only touch the text up to (not including) the paragraph
    #     mark so the edit is a clean in-place replacement.
    # ---------------------------------------------------------------
    $descText = "Want to play Clone Bonus slot for free? Read our comprehensive review, ratings, and pros and cons. Learn how to win big and make your free play today!"
    $bodyRange = $d.Range($imagePara.Range.Start, $imagePara.Range.End - 1)
    $bodyRange.Text = $descText

    # ---------------------------------------------------------------
    # 2b. Insert a brand-new paragraph right before it that repeats the
    #     page title in bold.
    # ---------------------------------------------------------------
    $imagePara = $d.Paragraphs.Item($imageIndex)
    $imagePara.Range.InsertParagraphBefore()

    $titlePara = $d.Paragraphs.Item($imageIndex)
    $titleText = "Play Clone Bonus Slot for Free - Review & Ratings 2021"
    $titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
    $titleRange.Text = $titleText
    $titleRange.Font.Bold = $true
    $titleRange.Font.Italic = $false
}
